$wb = $excel.ActiveWorkbook

# ---- 1. Update the summary ("总计") sheet: insert a new row for 2022-Q3 ----
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 4.78
# The row-insert kept each shifted row's original index (col A) attached to its
# content, but column A is really just a 0-based sequence number of the row -
# renumber rows 3-9 (1..7) so the index stays sequential after the insert.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---- 2. Insert a new worksheet "2022-Q3" before "2022-Q2" ----
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q3"

# ---- 3. Populate header row ----
$newSheet.Range("B1").Value = "'基金代码"
$newSheet.Range("C1").Value = "'基金名称"
$newSheet.Range("D1").Value = "'基金规模"
$newSheet.Range("E1").Value = "'股票总仓位"
$newSheet.Range("F1").Value = "'仓位占比"
$newSheet.Range("G1").Value = "'持有市值(亿元)"
$newSheet.Range("H1").Value = "'仓位排名"

# ---- 4. Populate data rows (A & H numeric, B-G forced text) ----
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'008099"
$newSheet.Range("C2").Value = "'广发价值领先混合A"
$newSheet.Range("D2").Value = "'58.59"
$newSheet.Range("E2").Value = "'93.49"
$newSheet.Range("F2").Value = "'5.13"
$newSheet.Range("G2").Value = "'3.0057"
$newSheet.Range("H2").Value = 5
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'012420"
$newSheet.Range("C3").Value = "'广发价值领先混合C"
$newSheet.Range("D3").Value = "'7.55"
$newSheet.Range("E3").Value = "'93.49"
$newSheet.Range("F3").Value = "'5.13"
$newSheet.Range("G3").Value = "'0.3873"
$newSheet.Range("H3").Value = 5
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'012528"
$newSheet.Range("C4").Value = "'广发鑫睿一年持有期混合A"
$newSheet.Range("D4").Value = "'7.26"
$newSheet.Range("E4").Value = "'92.69"
$newSheet.Range("F4").Value = "'5.10"
$newSheet.Range("G4").Value = "'0.3703"
$newSheet.Range("H4").Value = 8
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'013607"
$newSheet.Range("C5").Value = "'广发睿恒进取一年持有期混合A"
$newSheet.Range("D5").Value = "'6.95"
$newSheet.Range("E5").Value = "'91.91"
$newSheet.Range("F5").Value = "'5.14"
$newSheet.Range("G5").Value = "'0.3572"
$newSheet.Range("H5").Value = 5
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'014734"
$newSheet.Range("C6").Value = "'广发睿合混合A"
$newSheet.Range("D6").Value = "'5.96"
$newSheet.Range("E6").Value = "'86.96"
$newSheet.Range("F6").Value = "'5.29"
$newSheet.Range("G6").Value = "'0.3153"
$newSheet.Range("H6").Value = 6
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'012529"
$newSheet.Range("C7").Value = "'广发鑫睿一年持有期混合C"
$newSheet.Range("D7").Value = "'4.74"
$newSheet.Range("E7").Value = "'92.69"
$newSheet.Range("F7").Value = "'5.10"
$newSheet.Range("G7").Value = "'0.2417"
$newSheet.Range("H7").Value = 8
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'014735"
$newSheet.Range("C8").Value = "'广发睿合混合C"
$newSheet.Range("D8").Value = "'1.47"
$newSheet.Range("E8").Value = "'86.96"
$newSheet.Range("F8").Value = "'5.29"
$newSheet.Range("G8").Value = "'0.0778"
$newSheet.Range("H8").Value = 6
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'013608"
$newSheet.Range("C9").Value = "'广发睿恒进取一年持有期混合C"
$newSheet.Range("D9").Value = "'0.47"
$newSheet.Range("E9").Value = "'91.91"
$newSheet.Range("F9").Value = "'5.14"
$newSheet.Range("G9").Value = "'0.0242"
$newSheet.Range("H9").Value = 5

# ---- 5. Re-apply the canonical formatting (bold+border header row, bold+border column A) ----
# Use the sibling "2022-Q2" sheet (now shifted one position to the right) as the format template,
# since it already carries the exact same header/column-A styling as the new sheet should have.
$tmpl = $wb.Worksheets.Item("2022-Q2")
$tmpl.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$tmpl.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)

$newSheet.Range("A1").Select()
